# Working on making launch_intent results look nice.
# Add two new parameter columns (param.format, param.q) to the "queries" sheet header,
# and add a new "odk_values" query row (row 5) describing an odkquery:// lookup.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("queries")

# New header columns for query parameters
$ws.Range("D1").Value = "param.format"
$ws.Range("E1").Value = "param.q"

# New query row describing an ODK values lookup
$ws.Range("A5").Value = "odk_values"
$ws.Range("B5").Value = "odkquery://table_id/elementKey1/elementKey5/?selection=encodeURIComponent('elementKey2=? and elementKey3>5')&selectionArgs=encodeURIComponent(JSON.stringify([data('state')])"
